# Auto-generated script applying scheduled market-price/profit updates
# to the Ixion_Profits workbook (columns H-N) across sheets ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value2 = 249.86667
$ws.Range("I5").Value2 = 81.125
$ws.Range("J5").Value2 = 442.7143
$ws.Range("K5").Value2 = 81.125
$ws.Range("L5").Value2 = 442.7143
$ws.Range("M5").Value2 = 33.875
$ws.Range("N5").Value2 = -672.7143

$ws.Range("H11").Value2 = 515.13336
$ws.Range("I11").Value2 = 515.13336
$ws.Range("K11").Value2 = 515.13336
$ws.Range("M11").Value2 = -375.13336

$ws.Range("H20").Value2 = 3943.8
$ws.Range("I20").Value2 = 2874.75
$ws.Range("J20").Value2 = 8220
$ws.Range("K20").Value2 = 2874.75
$ws.Range("L20").Value2 = 8220
$ws.Range("M20").Value2 = -2644.75
$ws.Range("N20").Value2 = -8680

$ws.Range("H34").Value2 = 4907.7
$ws.Range("I34").Value2 = 2609.625
$ws.Range("K34").Value2 = 2609.625
$ws.Range("M34").Value2 = -2406.625

$ws.Range("H35").Value2 = 3943.8
$ws.Range("I35").Value2 = 2874.75
$ws.Range("J35").Value2 = 8220
$ws.Range("K35").Value2 = 2874.75
$ws.Range("L35").Value2 = 8220
$ws.Range("M35").Value2 = -2495.75
$ws.Range("N35").Value2 = -8978

$ws.Range("H36").Value2 = 4907.7
$ws.Range("I36").Value2 = 2609.625
$ws.Range("K36").Value2 = 2609.625
$ws.Range("M36").Value2 = -1894.625

$ws.Range("H127").Value2 = 1176
$ws.Range("I127").Value2 = 1029
$ws.Range("J127").Value2 = 1617
$ws.Range("K127").Value2 = 3087
$ws.Range("L127").Value2 = 4851
$ws.Range("M127").Value2 = 1873
$ws.Range("N127").Value2 = -14771

$ws.Range("H138").Value2 = 2052.6428
$ws.Range("I138").Value2 = 1250.3334
$ws.Range("J138").Value2 = 2582.9832
$ws.Range("K138").Value2 = 3751.0002
$ws.Range("L138").Value2 = 7748.9496
$ws.Range("M138").Value2 = 1388.9998
$ws.Range("N138").Value2 = -18028.9496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value2 = 4914.2856
$ws.Range("I3").Value2 = 2225
$ws.Range("J3").Value2 = 8500
$ws.Range("K3").Value2 = 2225
$ws.Range("L3").Value2 = 8500
$ws.Range("M3").Value2 = -2110
$ws.Range("N3").Value2 = -8730

$ws.Range("H5").Value2 = 170
$ws.Range("I5").Value2 = 140
$ws.Range("J5").Value2 = 290
$ws.Range("K5").Value2 = 140
$ws.Range("L5").Value2 = 290
$ws.Range("M5").Value2 = -28
$ws.Range("N5").Value2 = -514

$ws.Range("H8").Value2 = 8000
$ws.Range("J8").Value2 = 0
$ws.Range("L8").Value2 = 0
$ws.Range("N8").ClearContents()

$ws.Range("H45").Value2 = 7163.706
$ws.Range("I45").Value2 = 7561.4375
$ws.Range("J45").Value2 = 800
$ws.Range("K45").Value2 = 7561.4375
$ws.Range("L45").Value2 = 800
$ws.Range("M45").Value2 = -7184.4375
$ws.Range("N45").Value2 = -1554

$ws.Range("H61").Value2 = 179528.19
$ws.Range("I61").Value2 = 4659.6855
$ws.Range("J61").Value2 = 457728.1
$ws.Range("K61").Value2 = 4659.6855
$ws.Range("L61").Value2 = 457728.1
$ws.Range("M61").Value2 = -4447.6855
$ws.Range("N61").Value2 = -458152.1

$ws.Range("H136").Value2 = 179528.19
$ws.Range("I136").Value2 = 4659.6855
$ws.Range("J136").Value2 = 457728.1
$ws.Range("K136").Value2 = 13979.0565
$ws.Range("L136").Value2 = 1373184.3
$ws.Range("M136").Value2 = -11429.0565
$ws.Range("N136").Value2 = -1378284.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 170
$ws.Range("I4").Value2 = 140
$ws.Range("J4").Value2 = 290
$ws.Range("K4").Value2 = 140
$ws.Range("L4").Value2 = 290
$ws.Range("M4").Value2 = -25
$ws.Range("N4").Value2 = -520

$ws.Range("H8").Value2 = 1477.25
$ws.Range("I8").Value2 = 1401.8572
$ws.Range("J8").Value2 = 2005
$ws.Range("K8").Value2 = 1401.8572
$ws.Range("L8").Value2 = 2005
$ws.Range("M8").Value2 = -1261.8572
$ws.Range("N8").Value2 = -2285

$ws.Range("H22").Value2 = 382.5
$ws.Range("I22").Value2 = 338.7143
$ws.Range("K22").Value2 = 338.7143
$ws.Range("M22").Value2 = -165.7143

$ws.Range("H29").Value2 = 4311.6665
$ws.Range("I29").Value2 = 4311.6665
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 4311.6665
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = -4022.6665
$ws.Range("N29").ClearContents()

$ws.Range("H86").Value2 = 1662.7727
$ws.Range("I86").Value2 = 1599.05
$ws.Range("J86").Value2 = 2300
$ws.Range("K86").Value2 = 1599.05
$ws.Range("L86").Value2 = 2300
$ws.Range("M86").Value2 = -476.05
$ws.Range("N86").Value2 = -4546

$ws.Range("H89").Value2 = 1662.7727
$ws.Range("I89").Value2 = 1599.05
$ws.Range("J89").Value2 = 2300
$ws.Range("K89").Value2 = 7995.25
$ws.Range("L89").Value2 = 11500
$ws.Range("M89").Value2 = -2379.25
$ws.Range("N89").Value2 = -22732

$ws.Range("H94").Value2 = 2071.35
$ws.Range("I94").Value2 = 1617.3334
$ws.Range("J94").Value2 = 2752.375
$ws.Range("K94").Value2 = 1617.3334
$ws.Range("L94").Value2 = 2752.375
$ws.Range("M94").Value2 = -1166.3334
$ws.Range("N94").Value2 = -3654.375

$ws.Range("H105").Value2 = 2537.5
$ws.Range("I105").Value2 = 2600
$ws.Range("K105").Value2 = 2600
$ws.Range("M105").Value2 = -853

$ws.Range("H107").Value2 = 988.45
$ws.Range("I107").Value2 = 954.2143
$ws.Range("J107").Value2 = 1068.3334
$ws.Range("K107").Value2 = 954.2143
$ws.Range("L107").Value2 = 1068.3334
$ws.Range("M107").Value2 = 965.7857
$ws.Range("N107").Value2 = -4908.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value2 = 6000
$ws.Range("I93").Value2 = 6000
$ws.Range("K93").Value2 = 6000
$ws.Range("M93").Value2 = -4128

$ws.Range("H94").Value2 = 7031.778
$ws.Range("I94").Value2 = 3998.3333
$ws.Range("J94").Value2 = 8548.5
$ws.Range("K94").Value2 = 3998.3333
$ws.Range("L94").Value2 = 8548.5
$ws.Range("M94").Value2 = -3547.3333
$ws.Range("N94").Value2 = -9450.5

$ws.Range("H99").Value2 = 12496
$ws.Range("I99").Value2 = 15195
$ws.Range("K99").Value2 = 15195
$ws.Range("M99").Value2 = -13697

$ws.Range("H126").Value2 = 12496
$ws.Range("I126").Value2 = 15195
$ws.Range("K126").Value2 = 45585
$ws.Range("M126").Value2 = -43115

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 5551.04
$ws.Range("I5").Value2 = 8714.666999999999
$ws.Range("J5").Value2 = 2630.7693
$ws.Range("K5").Value2 = 26144.001
$ws.Range("L5").Value2 = 7892.3079
$ws.Range("M5").Value2 = -26032.001
$ws.Range("N5").Value2 = -8116.3079

$ws.Range("H117").Value2 = 25651282
$ws.Range("J117").Value2 = 37040570
$ws.Range("L117").Value2 = 111121710
$ws.Range("N117").Value2 = -111128594

$ws.Range("H121").Value2 = 10838.9375
$ws.Range("I121").Value2 = 297.83334
$ws.Range("J121").Value2 = 17163.6
$ws.Range("K121").Value2 = 893.5000200000001
$ws.Range("L121").Value2 = 51490.8
$ws.Range("M121").Value2 = 416.4999799999999
$ws.Range("N121").Value2 = -54110.8

$ws.Range("H135").Value2 = 5551.04
$ws.Range("I135").Value2 = 8714.666999999999
$ws.Range("J135").Value2 = 2630.7693
$ws.Range("K135").Value2 = 78432.003
$ws.Range("L135").Value2 = 23676.9237
$ws.Range("M135").Value2 = -75897.003
$ws.Range("N135").Value2 = -28746.9237

$ws.Range("H140").Value2 = 4624.278
$ws.Range("I140").Value2 = 4624.278
$ws.Range("K140").Value2 = 13872.834
$ws.Range("M140").Value2 = -8692.834000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 278.26666
$ws.Range("I2").Value2 = 224
$ws.Range("J2").Value2 = 305.4
$ws.Range("K2").Value2 = 224
$ws.Range("L2").Value2 = 305.4
$ws.Range("M2").Value2 = -111
$ws.Range("N2").Value2 = -531.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 20843602
$ws.Range("I132").Value2 = 22232774
$ws.Range("K132").Value2 = 66698322
$ws.Range("M132").Value2 = -66695792

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 64832.6
$ws.Range("I62").Value2 = 93649
$ws.Range("J62").Value2 = 7199.8
$ws.Range("K62").Value2 = 93649
$ws.Range("L62").Value2 = 7199.8
$ws.Range("M62").Value2 = -93025
$ws.Range("N62").Value2 = -8447.799999999999

$ws.Range("H65").Value2 = 64832.6
$ws.Range("I65").Value2 = 93649
$ws.Range("J65").Value2 = 7199.8
$ws.Range("K65").Value2 = 468245
$ws.Range("L65").Value2 = 35999
$ws.Range("M65").Value2 = -465125
$ws.Range("N65").Value2 = -42239

$ws.Range("H70").Value2 = 29600
$ws.Range("J70").Value2 = 30105
$ws.Range("L70").Value2 = 30105
$ws.Range("N70").Value2 = -30735

$ws.Range("H73").Value2 = 29600
$ws.Range("J73").Value2 = 30105
$ws.Range("L73").Value2 = 30105
$ws.Range("N73").Value2 = -32289

$ws.Range("H100").Value2 = 365
$ws.Range("I100").Value2 = 365
$ws.Range("K100").Value2 = 730
$ws.Range("M100").Value2 = -189

$ws.Range("H132").Value2 = 1627.1111
$ws.Range("I132").Value2 = 1322.2916
$ws.Range("J132").Value2 = 2236.75
$ws.Range("K132").Value2 = 3966.8748
$ws.Range("L132").Value2 = 6710.25
$ws.Range("M132").Value2 = -1436.8748
$ws.Range("N132").Value2 = -11770.25

